# Add an example dataset (one new paper + its associated measurement rows)
# to the "Data" and "Papers" sheets, as a reference example for Jerry.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("Data")
$papers = $wb.Worksheets.Item("Papers")

# --- Seed the new shared strings in the same order the original commit
# --- introduced them, so the shared-string table lines up: 130-140, M206,
# --- "218-300 ", the paper title, then the paper link.
$data.Range("C16").Value = "130-140"
$data.Range("D2").Value = "M206"
$data.Range("G2").Value = "218-300 "
$papers.Range("B2").Value = "Greenhouse gas emissions altered by the introduction of a year-long fallow to continuous rice systems"
$papers.Range("C2").Value = "https://doi.org/10.1002/jeq2.70055"

# --- Papers sheet: register the new paper as record #1.
$papers.Range("A2").Value = 1

# --- Data sheet: the measurement rows belonging to paper #1 (variety M206,
# --- year 2021) pulled from that paper.
$data.Range("A2").Value = 1
$data.Range("B2").Value = 0
$data.Range("C2").Value = 0
$data.Range("E2").Value = 2021

$data.Range("A3").Value = 1
$data.Range("B3").Value = 2
$data.Range("C3").Value = 7
$data.Range("D3").Value = "M206"
$data.Range("E3").Value = 2021
$data.Range("G3").Value = 250

$data.Range("A4").Value = 1
$data.Range("B4").Value = 10
$data.Range("C4").Value = 14
$data.Range("D4").Value = "M206"
$data.Range("E4").Value = 2021
$data.Range("G4").Value = 250

$data.Range("A5").Value = 1
$data.Range("D5").Value = "M206"
$data.Range("E5").Value = 2021
$data.Range("G5").Value = 250

$data.Range("A6").Value = 1
$data.Range("D6").Value = "M206"
$data.Range("E6").Value = 2021
$data.Range("G6").Value = 250

$data.Range("A7").Value = 1
$data.Range("D7").Value = "M206"
$data.Range("E7").Value = 2021
$data.Range("G7").Value = 250

$data.Range("A8").Value = 1
$data.Range("D8").Value = "M206"
$data.Range("E8").Value = 2021
$data.Range("G8").Value = 250

$data.Range("B9").Value = 6000
$data.Range("D9").Value = "M206"
$data.Range("E9").Value = 2021
$data.Range("G9").Value = 250

$data.Range("D10").Value = "M206"
$data.Range("E10").Value = 2021
$data.Range("G10").Value = 250

$data.Range("D11").Value = "M206"
$data.Range("E11").Value = 2021
$data.Range("G11").Value = 250

$data.Range("D12").Value = "M206"
$data.Range("E12").Value = 2021
$data.Range("G12").Value = 250

$data.Range("D13").Value = "M206"
$data.Range("E13").Value = 2021
$data.Range("G13").Value = 250

$data.Range("D14").Value = "M206"
$data.Range("E14").Value = 2021
$data.Range("G14").Value = 250

$data.Range("D15").Value = "M206"
$data.Range("E15").Value = 2021
$data.Range("G15").Value = 250

$data.Range("D16").Value = "M206"
$data.Range("E16").Value = 2021
$data.Range("G16").Value = 250

$data.Range("E18").Value = 2022
